$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 00:16"

$countryData = @(
    @(4, "China", 81171, 0, 73159, 4735, 1573, 0, 3277),
    @(5, "Italia", 69176, 5249, 8326, 54030, 3393, 743, 6820),
    @(6, "Estados Unidos", 53260, 9526, 370, 52201, 1175, 136, 689),
    @(7, "España", 39885, 4749, 3794, 33283, 2355, 497, 2808),
    @(8, "Alemania", 32986, 3930, 3243, 29586, 23, 34, 157),
    @(9, "Iran", 24811, 1762, 8913, 13964, 0, 122, 1934),
    @(10, "Francia", 22304, 2448, 3281, 17923, 2516, 240, 1100),
    @(11, "Suiza", 9877, 1082, 131, 9624, 141, 2, 122),
    @(12, "Corea del Sur", 9037, 76, 3507, 5410, 59, 9, 120),
    @(13, "Reino Unido", 8077, 1427, 135, 7520, 20, 87, 422),
    @(14, "Paises Bajos", 5560, 811, 2, 5282, 435, 63, 276),
    @(15, "Austria", 5283, 809, 9, 5246, 22, 7, 28),
    @(16, "Belgica", 4269, 526, 461, 3686, 381, 34, 122),
    @(17, "Noruega", 2863, 238, 6, 2845, 44, 2, 12),
    @(18, "Canada", 2792, 701, 112, 2654, 1, 2, 26),
    @(19, "Portugal", 2362, 302, 22, 2307, 48, 10, 33),
    @(20, "Suecia", 2286, 240, 16, 2234, 136, 9, 36),
    @(21, "Brasil", 2247, 323, 2, 2199, 18, 12, 46),
    @(22, "Australia", 2144, 257, 118, 2018, 11, 1, 8),
    @(23, "Israel", 1930, 488, 53, 1874, 34, 2, 3),
    @(24, "Turquia", 1872, 343, 0, 1828, 0, 7, 44),
    @(25, "Malasia", 1624, 106, 183, 1425, 64, 2, 16),
    @(26, "Dinamarca", 1591, 131, 1, 1558, 69, 8, 32),
    @(27, "Chequia", 1394, 158, 10, 1381, 19, 2, 3),
    @(28, "Irlanda", 1329, 204, 5, 1317, 29, 1, 7),
    @(29, "Japon", 1193, 65, 285, 865, 54, 1, 43),
    @(30, "Luxemburgo", 1099, 224, 6, 1085, 3, 0, 8),
    @(31, "Ecuador", 1049, 68, 3, 1019, 2, 9, 27),
    @(32, "Pakistan", 972, 97, 18, 947, 5, 1, 7),
    @(33, "Chile", 922, 176, 17, 903, 7, 0, 2),
    @(34, "Polonia", 901, 152, 1, 890, 3, 2, 10),
    @(35, "Tailandia", 827, 106, 52, 771, 7, 3, 4),
    @(36, "Rumania", 794, 218, 79, 704, 15, 4, 11),
    @(37, "Finlandia", 792, 92, 10, 781, 11, 0, 1),
    @(38, "Arabia Saudita", 767, 205, 28, 738, 0, 1, 1),
    @(39, "Grecia", 743, 48, 29, 694, 35, 3, 20),
    @(40, "Crucero", 712, 0, 587, 115, 15, 2, 10),
    @(41, "Indonesia", 686, 107, 30, 601, 0, 6, 55),
    @(42, "Islandia", 648, 60, 51, 595, 13, 1, 2),
    @(43, "Singapur", 558, 49, 156, 400, 14, 0, 2),
    @(44, "Sudafrica", 554, 152, 4, 550, 2, 0, 0),
    @(45, "Filipinas", 552, 90, 20, 497, 1, 2, 35),
    @(46, "India", 536, 37, 40, 486, 0, 0, 10),
    @(47, "Catar", 526, 25, 41, 485, 6, 0, 0),
    @(48, "Rusia", 495, 57, 22, 472, 8, 0, 1),
    @(49, "Eslovenia", 480, 38, 3, 473, 12, 1, 4),
    @(50, "Peru", 416, 21, 1, 408, 9, 2, 7),
    @(51, "Egipto", 402, 36, 80, 302, 0, 1, 20),
    @(52, "Barein", 392, 15, 177, 212, 2, 1, 3),
    @(53, "Hong Kong", 386, 29, 102, 280, 4, 0, 4),
    @(54, "Croacia", 382, 67, 5, 376, 6, 0, 1),
    @(55, "Colombia", 378, 101, 6, 369, 0, 0, 3),
    @(56, "Estonia", 369, 17, 7, 362, 4, 0, 0),
    @(57, "Mexico", 367, 51, 4, 359, 1, 1, 4),
    @(58, "Panama", 345, 0, 1, 338, 33, 0, 6),
    @(59, "Libano", 318, 51, 8, 306, 4, 0, 4),
    @(60, "Irak", 316, 50, 75, 214, 0, 4, 27),
    @(61, "Republica Dominicana", 312, 67, 3, 303, 0, 3, 6),
    @(62, "Serbia", 303, 54, 15, 285, 21, 0, 3),
    @(63, "Argentina", 301, 0, 52, 243, 0, 2, 6),
    @(64, "Argelia", 264, 34, 24, 221, 0, 2, 19),
    @(65, "Armenia", 249, 14, 14, 235, 6, 0, 0),
    @(66, "Emiratos Arabes Unidos", 248, 50, 45, 201, 2, 0, 2),
    @(67, "Bulgaria", 218, 17, 3, 212, 8, 0, 3),
    @(68, "Taiwan", 216, 21, 29, 185, 0, 0, 2),
    @(69, "Lituania", 209, 30, 1, 206, 1, 1, 2),
    @(70, "Eslovaquia", 204, 18, 7, 197, 2, 0, 0),
    @(71, "Letonia", 197, 17, 1, 196, 0, 0, 0),
    @(72, "Kuwait", 191, 2, 39, 152, 5, 0, 0),
    @(73, "San Marino", 187, 0, 4, 162, 12, 1, 21),
    @(74, "Hungria", 187, 20, 21, 157, 6, 1, 9),
    @(75, "Costa Rica", 177, 19, 2, 173, 4, 0, 2),
    @(76, "Marruecos", 170, 27, 6, 159, 1, 1, 5),
    @(77, "Bosnia y Herzegovina", 166, 30, 2, 161, 1, 2, 3),
    @(78, "Principado de Andorra", 164, 31, 1, 162, 7, 0, 1),
    @(79, "Uruguay", 162, 0, 0, 162, 3, 0, 0),
    @(80, "Nueva Zelanda", 155, 53, 12, 143, 0, 0, 0),
    @(81, "Jordania", 154, 27, 1, 153, 0, 0, 0),
    @(82, "Republica de Macedonia", 148, 12, 1, 145, 1, 0, 2),
    @(83, "Vietnam", 134, 11, 17, 117, 3, 0, 0),
    @(84, "Moldavia", 125, 16, 2, 122, 10, 0, 1),
    @(85, "Republica de Chipre", 124, 8, 3, 118, 3, 2, 3),
    @(86, "Albania", 123, 19, 10, 108, 2, 1, 5),
    @(87, "Islas Feroe", 122, 4, 33, 89, 0, 0, 0),
    @(88, "Tunez", 114, 25, 1, 109, 11, 1, 4),
    @(89, "Burkina Faso", 114, 15, 7, 103, 0, 0, 4),
    @(90, "Malta", 110, 3, 2, 108, 1, 0, 0),
    @(91, "Brunei", 104, 13, 2, 102, 2, 0, 0),
    @(92, "Sri Lanka", 102, 5, 2, 100, 2, 0, 0),
    @(93, "Ucrania", 97, 24, 1, 93, 0, 0, 3),
    @(94, "Reunion", 94, 23, 1, 93, 0, 0, 0),
    @(95, "Camboya", 91, 4, 4, 87, 1, 0, 0),
    @(96, "Azerbaiyan", 87, 15, 10, 76, 6, 0, 1),
    @(97, "Senegal", 86, 7, 8, 78, 0, 0, 0),
    @(98, "Venezuela", 84, 0, 15, 69, 2, 0, 0),
    @(99, "Oman", 84, 18, 17, 67, 0, 0, 0),
    @(100, "Bielorrusia", 81, 0, 22, 59, 0, 0, 0),
    @(101, "Afganistan", 74, 34, 1, 72, 0, 0, 1),
    @(102, "Costa de Marfil", 73, 48, 2, 71, 0, 0, 0),
    @(103, "Kazajistan", 72, 10, 0, 72, 0, 0, 0),
    @(104, "Georgia", 70, 9, 9, 61, 1, 0, 0),
    @(105, "Camerun", 66, 10, 2, 64, 0, 0, 0),
    @(106, "Guadalupe", 62, 0, 0, 61, 4, 0, 1),
    @(107, "Estado de Palestina", 60, 1, 16, 44, 0, 0, 0),
    @(108, "Trinidad yTobago", 57, 6, 0, 57, 0, 0, 0),
    @(109, "Martinica", 57, 4, 0, 56, 7, 0, 1),
    @(110, "Ghana", 53, 26, 0, 51, 0, 0, 2),
    @(111, "Liechtenstein", 51, 0, 0, 51, 0, 0, 0),
    @(112, "Uzbekistan", 50, 4, 0, 50, 0, 0, 0),
    @(113, "Cuba", 48, 8, 1, 46, 2, 0, 1),
    @(114, "Montenegro", 47, 20, 0, 46, 0, 0, 1),
    @(115, "Consejo Danes para los Refugiados", 45, 9, 0, 43, 0, 0, 2),
    @(116, "Nigeria", 44, 4, 2, 41, 0, 0, 1),
    @(117, "Kirguistan", 42, 26, 0, 42, 0, 0, 0),
    @(118, "Mauricio", 42, 6, 0, 40, 1, 0, 2),
    @(119, "Ruanda", 40, 4, 0, 40, 0, 0, 0),
    @(120, "Puerto Rico", 39, 8, 1, 36, 0, 0, 2),
    @(121, "Banglades", 39, 6, 5, 30, 0, 1, 4),
    @(122, "Mayotte", 36, 12, 0, 36, 0, 0, 0),
    @(123, "Guam", 32, 3, 0, 31, 0, 0, 1),
    @(124, "Honduras", 30, 0, 0, 30, 0, 0, 0),
    @(125, "Bolivia", 29, 2, 0, 29, 0, 0, 0),
    @(126, "Macao", 28, 3, 10, 18, 0, 0, 0),
    @(127, "Paraguay", 27, 5, 0, 25, 1, 1, 2),
    @(128, "Polinesia Francesa", 25, 7, 0, 25, 0, 0, 0),
    @(129, "Kenia", 25, 9, 0, 25, 0, 0, 0),
    @(130, "Isla de Man", 23, 10, 0, 23, 0, 0, 0),
    @(131, "Monaco", 23, 0, 1, 22, 0, 0, 0),
    @(132, "Guayana Francesa", 23, 3, 6, 17, 0, 0, 0),
    @(133, "Guatemala", 21, 1, 0, 20, 0, 0, 1),
    @(134, "Jamaica", 21, 2, 2, 18, 0, 0, 1),
    @(135, "Togo", 20, 2, 1, 19, 0, 0, 0),
    @(136, "Barbados", 18, 1, 0, 18, 0, 0, 0),
    @(137, "Madagascar", 17, 5, 0, 17, 0, 0, 0),
    @(138, "Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
    @(139, "Gibraltar", 15, 0, 5, 10, 0, 0, 0),
    @(140, "Maldivas", 13, 0, 5, 8, 0, 0, 0),
    @(141, "Etiopia", 12, 1, 0, 12, 0, 0, 0),
    @(142, "Tanzania", 12, 0, 0, 12, 0, 0, 0),
    @(143, "Aruba", 12, 3, 1, 11, 0, 0, 0),
    @(144, "Mongolia", 10, 0, 0, 10, 0, 0, 0),
    @(145, "Nueva Caledonia", 10, 2, 0, 10, 0, 0, 0),
    @(146, "Guinea Ecuatorial", 9, 0, 0, 9, 0, 0, 0),
    @(147, "Uganda", 9, 0, 0, 9, 0, 0, 0),
    @(148, "San Martin (Parte Francesa)", 8, 0, 0, 8, 0, 0, 0),
    @(149, "Surinam", 7, 2, 0, 7, 0, 0, 0),
    @(150, "Haiti", 7, 1, 0, 7, 0, 0, 0),
    @(151, "Seychelles", 7, 0, 0, 7, 0, 0, 0),
    @(152, "Namibia", 7, 3, 2, 5, 0, 0, 0),
    @(153, "Benin", 6, 0, 0, 6, 0, 0, 0),
    @(154, "Bermudas", 6, 0, 0, 6, 0, 0, 0),
    @(155, "Gabon", 6, 0, 0, 5, 0, 0, 1),
    @(156, "Curazao", 6, 2, 0, 5, 0, 0, 1),
    @(157, "Islas Caimanes", 6, 1, 0, 5, 0, 0, 1),
    @(158, "El Salvador", 5, 2, 0, 5, 0, 0, 0),
    @(159, "Guyana", 5, 0, 0, 4, 0, 0, 1),
    @(160, "Bahamas", 5, 1, 1, 4, 0, 0, 0),
    @(161, "Groenlandia", 5, 1, 2, 3, 0, 0, 0),
    @(162, "Guinea", 4, 0, 0, 4, 0, 0, 0),
    @(163, "Suazilandia", 4, 0, 0, 4, 0, 0, 0),
    @(164, "Santa Sede", 4, 3, 0, 4, 0, 0, 0),
    @(165, "Fiyi", 4, 1, 0, 4, 0, 0, 0),
    @(166, "Congo", 4, 0, 0, 4, 0, 0, 0),
    @(167, "Republica del Chad", 3, 1, 0, 3, 0, 0, 0),
    @(168, "Mozambique", 3, 2, 0, 3, 0, 0, 0),
    @(169, "Antigua y Barbuda", 3, 0, 0, 3, 0, 0, 0),
    @(170, "Liberia", 3, 0, 0, 3, 0, 0, 0),
    @(171, "Angola", 3, 0, 0, 3, 0, 0, 0),
    @(172, "Niger", 3, 0, 0, 3, 0, 0, 0),
    @(173, "Birmania", 3, 1, 0, 3, 0, 0, 0),
    @(174, "Santa Lucia", 3, 0, 0, 3, 0, 0, 0),
    @(175, "Zambia", 3, 0, 0, 3, 0, 0, 0),
    @(176, "Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0),
    @(177, "Republica de Yibuti", 3, 0, 0, 3, 0, 0, 0),
    @(178, "San Bartolome", 3, 0, 0, 3, 0, 0, 0),
    @(179, "Zimbabue", 3, 0, 0, 2, 0, 0, 1),
    @(180, "Cabo Verde", 3, 0, 0, 2, 0, 1, 1),
    @(181, "Gambia", 3, 1, 0, 2, 0, 0, 1),
    @(182, "Sudan", 3, 1, 0, 2, 0, 0, 1),
    @(183, "Laos", 2, 2, 0, 2, 0, 0, 0),
    @(184, "San Martin (Parte Holandesa)", 2, 0, 0, 2, 0, 0, 0),
    @(185, "Butan", 2, 0, 0, 2, 0, 0, 0),
    @(186, "Nicaragua", 2, 0, 0, 2, 0, 0, 0),
    @(187, "Mauritania", 2, 0, 0, 2, 0, 0, 0),
    @(188, "Dominica", 2, 0, 0, 2, 0, 0, 0),
    @(189, "Nepal", 2, 0, 1, 1, 0, 0, 0),
    @(190, "Belice", 1, 0, 0, 1, 0, 0, 0),
    @(191, "Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0),
    @(192, "Timor Oriental", 1, 0, 0, 1, 0, 0, 0),
    @(193, "Eritrea", 1, 0, 0, 1, 0, 0, 0),
    @(194, "Siria", 1, 0, 0, 1, 0, 0, 0),
    @(195, "Granada", 1, 0, 0, 1, 0, 0, 0),
    @(196, "Libia", 1, 1, 0, 1, 0, 0, 0),
    @(197, "Islas Turcas y Caicos", 1, 0, 0, 1, 0, 0, 0),
    @(198, "Montserrat", 1, 0, 0, 1, 0, 0, 0),
    @(199, "San Vicente y las Granadinas", 1, 0, 0, 1, 0, 0, 0),
    @(200, "Somalia", 1, 0, 0, 1, 0, 0, 0)
)

foreach ($row in $countryData) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}

Write-Output "applied $($countryData.Count) rows"
